$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: four new formatted (but empty) cells (done first so the new
#     style slots line up the same way Excel allocated them originally) ---
$ws.Range("F21").VerticalAlignment = -4108
$ws.Range("H21").Value = "'"
$ws.Range("H21").ClearContents()

# --- Row 10: now holds what used to be row 11's data (Clau porta / Rfdi / (Wemos)) ---
$ws.Range("B10").Value = "Clau porta"
$ws.Range("C10").Value = "Rfdi"
$ws.Range("D10").Value = "(Wemos)"

# --- Row 11: now holds what used to be row 10's data, relabelled as the REBEDOR group head ---
$ws.Range("A11").Value = "REBEDOR"
$ws.Range("B11").Value = "Timbre entrada"
$ws.Range("C11").Value = "Bronzidor"
$ws.Range("D11").Value = "'Placa relé + ESP-01"

# Row 11 takes on the "first row of group" borders (thin box around A:D, no bottom)
$ws.Range("A11:D11").Borders.Item(7).LineStyle = 1
$ws.Range("A11:D11").Borders.Item(8).LineStyle = 1
$ws.Range("A11:D11").Borders.Item(10).LineStyle = 1

# --- Row 12: unchanged content, but style normalised to the "inner row" border (was bottom-border row) ---
$ws.Range("A12:D12").Borders.Item(9).LineStyle = -4142

# --- Row 13: loses the REBEDOR label + its "first row of group" border ---
$ws.Range("A13").ClearContents()
$ws.Range("A13:D13").Borders.Item(8).LineStyle = -4142

# --- Row 14: D14 loses its "(Wemos)" value ---
$ws.Range("D14").ClearContents()
